# Levi's CMS Operation&Maintenance Request List_CHN.xlsx
# Apply the 2019/11/26 & 2019/11/28 ticket updates to the "课题管理表" sheet:
#   - Row 11 (Wuhan F2 material ticket) gets marked completed: status -> 确认\n完成,
#     the execution-result note corrected to the Wuhan F2 text, a completion date,
#     an executiver and a reviewer are filled in, and the row is restyled to the
#     "completed" (filled / shaded) look used by rows 6-10.
#   - Row 12, previously a blank template row, becomes a brand-new fault ticket
#     (Wuhan A1 screen white-screen issue) with the full set of field values and
#     the same "completed" styling / row height as row 11.
#   - Misc view-state bits (selected cell) are refreshed to reflect the new
#     working row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 11 - close out the existing "Wuhan F2 material" ticket
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "确认`n完成"
$ws.Range("Q11").Value = "【2019/11/26】`n已为客户提取武汉店F2素材"
$ws.Range("R11").Value = 43795
$ws.Range("S11").Value = "李琳"
$ws.Range("T11").Value = "郭文博"

# ---------------------------------------------------------------------------
# 2) Row 12 - brand-new fault ticket (previously an empty template row)
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "提出"
$ws.Range("C12").Value = "故障对应"
$ws.Range("D12").Value = "0020028831"
$ws.Range("E12").Value = "WUHAN BEACON STREET"
$ws.Range("F12").Value = "China"
$ws.Range("G12").Value = "Wuhan"
$ws.Range("H12").Value = "154118"
$ws.Range("I12").Value = "2019/11/28"
$ws.Range("J12").Value = "2019/11/28"
$ws.Range("K12").Value = "H"
$ws.Range("L12").Value = "M"
$ws.Range("M12").Value = "S"
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = "谭志明"
$ws.Range("P12").Value = "【2019/11/28】`n武汉店A1屏幕出现白屏"
$ws.Range("Q12").Value = "【2019/11/28】`n远程重启设备后正常播放，判断为长时间未断电设备过热引起"
$ws.Range("R12").Value = 43797
$ws.Range("S12").Value = "刘琪"
$ws.Range("T12").Value = "郭文博"

# ---------------------------------------------------------------------------
# 3) Restyle rows 11 & 12 like the other completed rows (6-10): copy the
#    number formats / fills / borders from row 10 (already ht=36, same
#    column layout) onto A11:V11 and A12:V12, then fix up the row height
#    for row 12 (row 11 already is 36).
# ---------------------------------------------------------------------------
$styleSource = $ws.Range("A10:V10")
$styleSource.Copy()
$ws.Range("A11:V11").PasteSpecial(-4122)
$ws.Range("A12:V12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(11).RowHeight = 36
$ws.Rows.Item(12).RowHeight = 36

# ---------------------------------------------------------------------------
# 4) Refresh the active selection to the new working cell
# ---------------------------------------------------------------------------
$ws.Range("M9").Select()
